$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.309.61"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "1.707.37"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "223.18"
$ws.Range("E5").Value = "  -2.85%  "
$ws.Range("D6").Value = "0.5306"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("D8").Value = "0.2656"
$ws.Range("E8").Value = "  -4.66%  "
$ws.Range("D9").Value = "0.06590"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").Value = "20.88"
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").Value = "0.07645"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "4.586"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "1.735.39"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "1.942.08"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "0.5738"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").Value = "0.0₅8191"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "67.58"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").Value = "27.292.20"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "216.32"
$ws.Range("E19").Value = "  -4.42%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "4.674"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("D23").Value = "5.978"
$ws.Range("E23").Value = "  -4.80%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "142.30"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").Value = "1.749"
$ws.Range("E26").Value = "  +6.88%  "
$ws.Range("D27").Value = "0.1216"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").Value = "7.254"
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("D29").Value = "16.30"
$ws.Range("E29").Value = "  -4.25%  "
$ws.Range("D30").Value = "0.05372"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "3.510"
$ws.Range("E32").Value = "  -5.68%  "
$ws.Range("E33").Value = "  -3.63%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "2.878"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").Value = "2.421"
$ws.Range("D37").Value = "0.9474"
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("D38").Value = "0.5855"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "0.01632"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D40").Value = "5.863"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "1.040.54"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "0.8405"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "101.01"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "1.848.72"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "1.007"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "8.088"
$ws.Range("E50").Value = "  -2.56%  "
$ws.Range("D51").Value = "0.06536"
$ws.Range("E51").Value = "  +10.32%  "
